# Actualización automática de tasas-transfi.xlsx
$wb = $excel.ActiveWorkbook

# --- Hoja1: update the "Conversión del día" note with the new Binance rates ---
$ws1 = $wb.Worksheets.Item("Hoja1")

$newText = "Conversión del día 💰" + [char]10 + `
"✅ Dólar paralelo: 68" + [char]10 + `
"" + [char]10 + `
"Binance" + [char]10 + `
"✅ 1000 Bs = 9.41 = 39257.4 pesos" + [char]10 + `
"✅ 39257.4 pesos = 9.35 = 956.43 Bs" + [char]10 + `
"" + [char]10 + `
"Promedio competencia" + [char]10 + `
"✅ Tasa pesos: 20" + [char]10 + `
"✅ Tasa Bs: 20" + [char]10 + `
"✅ % Ganancia: 20%"

$ws1.Range("A1").Value = $newText

# --- tasas: update the rate-table cells ---
$ws2 = $wb.Worksheets.Item("tasas")

$ws2.Range("N10").Value = 106.222
$ws2.Range("O10").Value = 4170
$ws2.Range("N12").Value = 4199
$ws2.Range("O12").Value = 102.3
